# AStar - Test on empty 24x24 map
# Clear every interior (non-border) cell of the 24x24 map grid back to 0,
# leaving the outer wall border (row 1, row 24, column A, column X) intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2:W23").Value = 0

# Update the sheet's view/selection to match the saved state after the edit:
# scrolled down a few rows and the interior play-area selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
[void]$ws.Range("B2:W23").Select()
